# Scheduled-runner data refresh: update market/profit figures across all class sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match the latest pulled pricing data.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(5, 8).Value = 636.5
$ws.Cells.Item(5, 9).Value = 667.25
$ws.Cells.Item(5, 10).Value = 575
$ws.Cells.Item(5, 11).Value = 667.25
$ws.Cells.Item(5, 12).Value = 575
$ws.Cells.Item(5, 13).Value = -552.25
$ws.Cells.Item(5, 14).Value = -805
$ws.Cells.Item(70, 8).Value = 2116.5264
$ws.Cells.Item(70, 9).Value = 1565.4
$ws.Cells.Item(70, 11).Value = 4696.200000000001
$ws.Cells.Item(70, 13).Value = -4426.200000000001
$ws.Cells.Item(73, 8).Value = 2116.5264
$ws.Cells.Item(73, 9).Value = 1565.4
$ws.Cells.Item(73, 11).Value = 4696.200000000001
$ws.Cells.Item(73, 13).Value = -3760.200000000001
$ws.Cells.Item(100, 8).Value = 6026.8
$ws.Cells.Item(100, 9).Value = 2320.5557
$ws.Cells.Item(100, 10).Value = 11586.167
$ws.Cells.Item(100, 11).Value = 2320.5557
$ws.Cells.Item(100, 12).Value = 11586.167
$ws.Cells.Item(100, 13).Value = -1779.5557
$ws.Cells.Item(100, 14).Value = -12668.167
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(116, 8).Value = 199204.77
$ws.Cells.Item(116, 10).Value = 447289.62
$ws.Cells.Item(116, 12).Value = 447289.62
$ws.Cells.Item(116, 14).Value = -454173.62
$ws.Cells.Item(125, 8).Value = 1537
$ws.Cells.Item(125, 9).Value = 1189
$ws.Cells.Item(125, 11).Value = 10701
$ws.Cells.Item(125, 13).Value = -8241
$ws.Cells.Item(132, 8).Value = 77208.39
$ws.Cells.Item(132, 9).Value = 84820.17999999999
$ws.Cells.Item(132, 10).Value = 6165
$ws.Cells.Item(132, 11).Value = 254460.54
$ws.Cells.Item(132, 12).Value = 18495
$ws.Cells.Item(132, 13).Value = -251930.54
$ws.Cells.Item(132, 14).Value = -23555
$ws.Cells.Item(137, 8).Value = 9858.333000000001
$ws.Cells.Item(137, 9).Value = 3489.5
$ws.Cells.Item(137, 11).Value = 10468.5
$ws.Cells.Item(137, 13).Value = -7918.5
$ws.Cells.Item(141, 8).Value = 2949.111
$ws.Cells.Item(141, 9).Value = 2949.111
$ws.Cells.Item(141, 11).Value = 8847.332999999999
$ws.Cells.Item(141, 13).Value = -3667.332999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(132, 8).Value = 742636.25
$ws.Cells.Item(132, 9).Value = 871051.4
$ws.Cells.Item(132, 10).Value = 4249.25
$ws.Cells.Item(132, 11).Value = 2613154.2
$ws.Cells.Item(132, 12).Value = 12747.75
$ws.Cells.Item(132, 13).Value = -2610624.2
$ws.Cells.Item(132, 14).Value = -17807.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(99, 8).Value = 19419.467
$ws.Cells.Item(99, 9).Value = 21735.576
$ws.Cells.Item(99, 10).Value = 4364.75
$ws.Cells.Item(99, 11).Value = 21735.576
$ws.Cells.Item(99, 12).Value = 4364.75
$ws.Cells.Item(99, 13).Value = -20237.576
$ws.Cells.Item(99, 14).Value = -7360.75
$ws.Cells.Item(132, 8).Value = 6000000
$ws.Cells.Item(132, 10).Value = 6000000
$ws.Cells.Item(132, 12).Value = 6000000
$ws.Cells.Item(132, 14).Value = -6010120

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 18801.291
$ws.Cells.Item(31, 9).Value = 8705.75
$ws.Cells.Item(31, 10).Value = 25177.422
$ws.Cells.Item(31, 11).Value = 8705.75
$ws.Cells.Item(31, 12).Value = 25177.422
$ws.Cells.Item(31, 13).Value = -8410.75
$ws.Cells.Item(31, 14).Value = -25767.422
$ws.Cells.Item(34, 8).Value = 18801.291
$ws.Cells.Item(34, 9).Value = 8705.75
$ws.Cells.Item(34, 10).Value = 25177.422
$ws.Cells.Item(34, 11).Value = 8705.75
$ws.Cells.Item(34, 12).Value = 25177.422
$ws.Cells.Item(34, 13).Value = -8503.75
$ws.Cells.Item(34, 14).Value = -25581.422
$ws.Cells.Item(58, 8).Value = 497089
$ws.Cells.Item(58, 9).Value = 728730.3
$ws.Cells.Item(58, 10).Value = 4851.25
$ws.Cells.Item(58, 11).Value = 728730.3
$ws.Cells.Item(58, 12).Value = 4851.25
$ws.Cells.Item(58, 13).Value = -728527.3
$ws.Cells.Item(58, 14).Value = -5257.25
$ws.Cells.Item(132, 8).Value = 202207.88
$ws.Cells.Item(132, 9).Value = 1731.579
$ws.Cells.Item(132, 11).Value = 5194.737
$ws.Cells.Item(132, 13).Value = -2664.737
$ws.Cells.Item(134, 8).Value = 2280.4
$ws.Cells.Item(134, 9).Value = 2000.7142
$ws.Cells.Item(134, 10).Value = 2933
$ws.Cells.Item(134, 11).Value = 6002.142599999999
$ws.Cells.Item(134, 12).Value = 8799
$ws.Cells.Item(134, 13).Value = -3467.142599999999
$ws.Cells.Item(134, 14).Value = -13869
$ws.Cells.Item(136, 8).Value = 497089
$ws.Cells.Item(136, 9).Value = 728730.3
$ws.Cells.Item(136, 10).Value = 4851.25
$ws.Cells.Item(136, 11).Value = 2186190.9
$ws.Cells.Item(136, 12).Value = 14553.75
$ws.Cells.Item(136, 13).Value = -2183640.9
$ws.Cells.Item(136, 14).Value = -19653.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(107, 8).Value = 705.4167
$ws.Cells.Item(107, 9).Value = 671.625
$ws.Cells.Item(107, 10).Value = 773
$ws.Cells.Item(107, 11).Value = 2014.875
$ws.Cells.Item(107, 12).Value = 2319
$ws.Cells.Item(107, 13).Value = -94.875
$ws.Cells.Item(107, 14).Value = -6159
$ws.Cells.Item(131, 8).Value = 12912.421
$ws.Cells.Item(131, 10).Value = 29548.375
$ws.Cells.Item(131, 12).Value = 88645.125
$ws.Cells.Item(131, 14).Value = -98725.125

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(102, 8).Value = 2950.3333
$ws.Cells.Item(102, 9).Value = 2009.3043
$ws.Cells.Item(102, 11).Value = 2009.3043
$ws.Cells.Item(102, 13).Value = -387.3043
$ws.Cells.Item(132, 8).Value = 5301.2173
$ws.Cells.Item(132, 9).Value = 3487.5
$ws.Cells.Item(132, 10).Value = 8122.5557
$ws.Cells.Item(132, 11).Value = 10462.5
$ws.Cells.Item(132, 12).Value = 24367.6671
$ws.Cells.Item(132, 13).Value = -7932.5
$ws.Cells.Item(132, 14).Value = -29427.6671
$ws.Cells.Item(136, 8).Value = 30581.375
$ws.Cells.Item(136, 10).Value = 30581.375
$ws.Cells.Item(136, 12).Value = 91744.125
$ws.Cells.Item(136, 14).Value = -96844.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(40, 8).Value = 2319.8147
$ws.Cells.Item(40, 9).Value = 2321.64
$ws.Cells.Item(40, 10).Value = 2297
$ws.Cells.Item(40, 11).Value = 2321.64
$ws.Cells.Item(40, 12).Value = 2297
$ws.Cells.Item(40, 13).Value = -2185.64
$ws.Cells.Item(40, 14).Value = -2569
$ws.Cells.Item(81, 8).Value = 30000
$ws.Cells.Item(81, 10).Value = 30000
$ws.Cells.Item(81, 12).Value = 30000
$ws.Cells.Item(81, 14).Value = -31996
$ws.Cells.Item(84, 8).Value = 30000
$ws.Cells.Item(84, 10).Value = 30000
$ws.Cells.Item(84, 12).Value = 90000
$ws.Cells.Item(84, 14).Value = -99984
$ws.Cells.Item(100, 8).Value = 11618.546
$ws.Cells.Item(100, 9).Value = 2755.5557
$ws.Cells.Item(100, 11).Value = 2755.5557
$ws.Cells.Item(100, 13).Value = -2214.5557
$ws.Cells.Item(122, 8).Value = 4370.485
$ws.Cells.Item(122, 9).Value = 4081.1155
$ws.Cells.Item(122, 11).Value = 12243.3465
$ws.Cells.Item(122, 13).Value = -9793.3465
$ws.Cells.Item(132, 8).Value = 21555.385
$ws.Cells.Item(132, 9).Value = 26522.1
$ws.Cells.Item(132, 11).Value = 79566.29999999999
$ws.Cells.Item(132, 13).Value = -77036.29999999999
$ws.Cells.Item(136, 8).Value = 2494.5278
$ws.Cells.Item(136, 9).Value = 1796.8214
$ws.Cells.Item(136, 11).Value = 5390.4642
$ws.Cells.Item(136, 13).Value = -2840.4642

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(81, 8).Value = 1590
$ws.Cells.Item(81, 9).Value = 1453.5
$ws.Cells.Item(81, 10).Value = 1999.5
$ws.Cells.Item(81, 11).Value = 2907
$ws.Cells.Item(81, 12).Value = 3999
$ws.Cells.Item(81, 13).Value = -1846
$ws.Cells.Item(81, 14).Value = -6121
$ws.Cells.Item(84, 8).Value = 1590
$ws.Cells.Item(84, 9).Value = 1453.5
$ws.Cells.Item(84, 10).Value = 1999.5
$ws.Cells.Item(84, 11).Value = 14535
$ws.Cells.Item(84, 12).Value = 19995
$ws.Cells.Item(84, 13).Value = -9231
$ws.Cells.Item(84, 14).Value = -30603
$ws.Cells.Item(132, 8).Value = 3625992.2
$ws.Cells.Item(132, 9).Value = 4169291
$ws.Cells.Item(132, 10).Value = 4001.6667
$ws.Cells.Item(132, 11).Value = 12507873
$ws.Cells.Item(132, 12).Value = 12005.0001
$ws.Cells.Item(132, 13).Value = -12505343
$ws.Cells.Item(132, 14).Value = -17065.0001
$ws.Cells.Item(136, 8).Value = 8876.966
$ws.Cells.Item(136, 9).Value = 9953.772000000001
$ws.Cells.Item(136, 11).Value = 29861.316
$ws.Cells.Item(136, 13).Value = -27311.316
